$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Rows 3 through 24: decrement "PERIOD TO EXPIRE" (column H) by 1
# and update "LAST UPDATE" (column I) from 03-Nov-2025 to 04-Nov-2025
for ($row = 3; $row -le 24; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)
    $iCell.NumberFormat = "@"
    $iCell.Value2 = "04-Nov-2025"
}
